$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 51 raw survey answers were corrected (all 5s now) ---
$ws.Range("B51").Value = 5
$ws.Range("C51").Value = 5
$ws.Range("D51").Value = 5
$ws.Range("E51").Value = 5
$ws.Range("F51").Value = 5
$ws.Range("G51").Value = 5

# That row is no longer flagged as a "Bad" entry -> restore the default
# (Normal) cell style.
$ws.Range("B51:G51").Style = "Normal"

# The "Bad" cell style is now unused anywhere in the workbook, so delete it
# outright (this also drops its font/fill from the style table).
$wb.Styles.Item("Bad").Delete()

# --- Window/view state: sheet scrolled back up and a different cell
# selected ---
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("H14").Select()

# Recalculate the dependent SUM / AVERAGE formulas in rows 56-57.
$excel.Calculate()
